$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, taken from the latest crypto price/volume
# pull (including the PEPE/PancakeSwap and RenderToken/OKB row re-ordering).
$updates = [ordered]@{
    'D2' = '58.886.60'
    'E2' = '  -0.85%  '
    'D3' = '2.525.76'
    'E3' = '  +3.23%  '
    'E4' = '  +0.08%  '
    'D5' = '536.84'
    'E5' = '  +0.42%  '
    'D6' = '143.93'
    'E6' = '  -2.51%  '
    'D7' = '0.999'
    'E7' = '  +0.21%  '
    'E8' = '  +0.38%  '
    'D9' = '2.524.83'
    'E9' = '  +2.49%  '
    'D10' = '0.0995'
    'E10' = '  +0.39%  '
    'D12' = '5.52'
    'E12' = '  +2.41%  '
    'D13' = '0.350'
    'E13' = '  -0.15%  '
    'D14' = '2.966.81'
    'E14' = '  +3.09%  '
    'D15' = '23.51'
    'E15' = '  -2.70%  '
    'D16' = '58.891.03'
    'E16' = '  -0.76%  '
    'E17' = '  +0.26%  '
    'D18' = '2.518.22'
    'E18' = '  +1.03%  '
    'E19' = '  +0.25%  '
    'E20' = '  -2.62%  '
    'D21' = '322.74'
    'E21' = '  -0.72%  '
    'D22' = '1.00'
    'E22' = '  +3.06%  '
    'E23' = '  +0.90%  '
    'D24' = '61.76'
    'E24' = '  +2.29%  '
    'D25' = '0.436'
    'E25' = '  -6.79%  '
    'E26' = '  +0.67%  '
    'D27' = '2.625.33'
    'E27' = '  +2.75%  '
    'D28' = '0.995'
    'E28' = '  +2.00%  '
    'E29' = '  -0.24%  '
    'D30' = '6.71'
    'E30' = '  -2.77%  '
    'B31' = 'PEPE'
    'C31' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D31' = '0.0₃0767'
    'E31' = '  -0.48%  '
    'B32' = 'PancakeSwap'
    'C32' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D32' = '1.79'
    'E32' = '  -1.52%  '
    'E33' = '  -9.18%  '
    'E34' = '  +0.16%  '
    'D35' = '158.17'
    'E35' = '  +0.79%  '
    'D36' = '1.43'
    'E36' = '  +6.07%  '
    'D37' = '18.58'
    'E37' = '  +1.40%  '
    'D38' = '4.34'
    'E38' = '  -4.65%  '
    'D39' = '1.61'
    'E39' = '  -7.49%  '
    'B40' = 'RenderToken'
    'C40' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D40' = '5.58'
    'E40' = '  -3.20%  '
    'B41' = 'OKB'
    'C41' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D41' = '36.45'
    'E41' = '  -0.88%  '
    'D42' = '295.62'
    'E42' = '  -5.41%  '
    'D43' = '3.64'
    'E43' = '  -1.89%  '
    'D44' = '0.811'
    'E44' = '  -5.09%  '
    'D45' = '0.998'
    'E45' = '  +0.25%  '
    'D46' = '0.604'
    'E46' = '  +3.63%  '
    'D47' = '10.77'
    'E47' = '  +0.65%  '
    'D48' = '124.37'
    'E48' = '  +3.91%  '
    'D49' = '0.0929'
    'E49' = '  -1.05%  '
    'D50' = '18.60'
    'E50' = '  +0.54%  '
    'D51' = '0.0512'
    'E51' = '  -2.00%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text format first so numeric-looking strings such as "536.84" or
    # "1.00" are preserved as text instead of being coerced into numbers.
    $range.NumberFormat = '@'
    $range.Value = $updates[$cellRef]
}
